# Add a new "release/1.0.2" row to the meta-sheet, marking dev/sit/uat as
# not applicable ("X") while keeping prod's value, mirroring the existing
# "release/1.0.0" row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = "release/1.0.2"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
